$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the three "... Value" headers to "... Value in GBP" (also updates the
# bound ListObject/table column names since they share the same cells).
$ws.Range("D1").Value = "Buy Value in GBP"
$ws.Range("G1").Value = "Sell Value in GBP"
$ws.Range("J1").Value = "Fee Value in GBP"

# Give the header row a thin white border (matches the table's "header row"
# style band) and re-apply the bold/white-on-black header formatting.
$headerRange = $ws.Range("A1:M1")
$headerRange.Borders.Color = 16777215
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 0

# Autofit columns A:M now that header text changed (matches bestFit columns
# in the saved file).
$ws.Range("A1:M7").Columns.AutoFit()

# Move the active selection (cosmetic, matches the saved view state).
$ws.Range("D16").Select()
